$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.64"
$ws.Range("E2").Value = "'-0.55%"

$ws.Range("D3").Value = "'43.45"
$ws.Range("E3").Value = "'4.11%"

$ws.Range("D4").Value = "'5.604"
$ws.Range("E4").Value = "'-1.04%"

$ws.Range("D5").Value = "'0.08204"
$ws.Range("E5").Value = "'-2.05%"

$ws.Range("D6").Value = "'8.769"
$ws.Range("E6").Value = "'-0.26%"

$ws.Range("D7").Value = "'4.365"
$ws.Range("E7").Value = "'-4.08%"

$ws.Range("D8").Value = "'1.884"
$ws.Range("E8").Value = "'-6.54%"

$ws.Range("E9").Value = "'-5.72%"

$ws.Range("D10").Value = "'0.9434"
$ws.Range("E10").Value = "'1.66%"

$ws.Range("D11").Value = "'0.1188"
$ws.Range("E11").Value = "'-8.48%"

$ws.Range("D12").Value = "'0.1910"
$ws.Range("E12").Value = "'-2.70%"

$ws.Range("D13").Value = "'0.09750"
$ws.Range("E13").Value = "'3.95%"

$ws.Range("D14").Value = "'0.04321"
$ws.Range("E14").Value = "'11.03%"

$ws.Range("D15").Value = "'0.1070"
$ws.Range("E15").Value = "'0.79%"

$ws.Range("D16").Value = "'0.001281"
$ws.Range("E16").Value = "'-1.55%"

$ws.Range("D17").Value = "'0.005940"
$ws.Range("E17").Value = "'-2.56%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.531"
$ws.Range("E18").Value = "'2.72%"

$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3536"
$ws.Range("E19").Value = "'-0.07%"

$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'8.741"
$ws.Range("E20").Value = "'9.22%"

$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1370"
$ws.Range("E21").Value = "'-0.31%"

$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2495"
$ws.Range("E22").Value = "'-4.52%"

$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04396"
$ws.Range("E23").Value = "'-0.59%"

$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001240"
$ws.Range("E24").Value = "'-1.52%"

$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").Value = "'0.004301"
$ws.Range("E25").Value = "'-1.94%"

$ws.Range("D26").Value = "'0.0001237"
$ws.Range("E26").Value = "'2.91%"

$ws.Range("D27").Value = "'0.0004010"
$ws.Range("E27").Value = "'31.69%"

$ws.Range("D39").Value = "'0.02750"
$ws.Range("E39").Value = "'-1.50%"

$ws.Range("D40").Value = "'0.05681"
$ws.Range("E40").Value = "'2.59%"

$ws.Range("D41").Value = "'0.007901"
$ws.Range("E41").Value = "'1.37%"

$ws.Range("D42").Value = "'0.009758"
$ws.Range("E42").Value = "'4.54%"

$ws.Range("E43").Value = "'-1.28%"

$ws.Range("E44").Value = "'-2.53%"

$ws.Range("D45").Value = "'0.01005"
$ws.Range("E45").Value = "'-9.44%"

$ws.Range("D46").Value = "'0.00007324"
$ws.Range("E46").Value = "'4.34%"

$ws.Range("D47").Value = "'0.00000000754"
$ws.Range("E47").Value = "'0.32%"

$ws.Range("D48").Value = "'0.003449"
$ws.Range("E48").Value = "'-2.24%"

$ws.Range("D49").Value = "'0.002282"
$ws.Range("E49").Value = "'-0.01%"

$ws.Range("D50").Value = "'0.00002111"
$ws.Range("E50").Value = "'0.32%"

$ws.Range("D51").Value = "'0.0002011"
$ws.Range("E51").Value = "'0.32%"
